# Remove Alarm from diagram, no longer used for autoscaling.
#
# The architecture diagram on slide 1 had an "Alarm" icon + label (next to
# "Logs") that is no longer relevant, and the bottom "Darktrace appliance"
# caption is renamed to "Darktrace Master Instance".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. Remove the "Alarm" label textbox and its bar-chart icon picture ---
# Walk backwards so deleting earlier items doesn't perturb later indices.
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s.Shapes.Item($i)

    $isAlarmLabel = ($sh.Name -eq "TextBox 16") -and $sh.HasTextFrame -and $sh.TextFrame.HasText -and ($sh.TextFrame.TextRange.Text -eq "Alarm")

    if ($isAlarmLabel) {
        # The small graphic that sits just above/with the "Alarm" caption
        # (the next shape in z-order, use the 1-based collection position).
        $iconIndex = $i + 1
        if ($iconIndex -le $s.Shapes.Count) {
            $icon = $s.Shapes.Item($iconIndex)
            if ($icon.Name -eq "Graphic 6") {
                $icon.Delete()
            }
        }
        $sh.Delete()
    }
}

# --- 2. Rename the "Darktrace appliance" caption to "Darktrace Master Instance" ---
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Rectangle 1" -and $sh.HasTextFrame -and $sh.TextFrame.HasText -and ($sh.TextFrame.TextRange.Text -eq "Darktrace appliance")) {
        $sh.TextFrame.TextRange.Text = "Darktrace Master Instance"
        # Shape uses spAutoFit; pin the resulting box to the exact size/pos
        # PowerPoint computes for the new caption (grows to fit on one line,
        # staying vertically anchored and centered on the same midpoint).
        $sh.Left = 465.1196850393701
        $sh.Top = 499.6776377952756
        $sh.Width = 159.6287401574803
        $sh.Height = 21.810944881889764
    }
}
